$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 262
$ws.Cells.Item(2, 3).Value = 4
$ws.Cells.Item(2, 4).Value = 606
$ws.Cells.Item(2, 5).Value = 634
$ws.Cells.Item(2, 6).Value = 7
$ws.Cells.Item(2, 7).Value = 1
$ws.Cells.Item(2, 8).Value = 1902
$ws.Cells.Item(2, 9).Value = 31044
$ws.Cells.Item(2, 10).Value = 9
$ws.Cells.Item(2, 11).Value = 17
$ws.Cells.Item(2, 12).Value = 909

$ws.Cells.Item(3, 2).Value = 271
$ws.Cells.Item(3, 3).Value = 0
$ws.Cells.Item(3, 4).Value = 0
$ws.Cells.Item(3, 5).Value = 0
$ws.Cells.Item(3, 6).Value = 0
$ws.Cells.Item(3, 7).Value = 0
$ws.Cells.Item(3, 8).Value = 0
$ws.Cells.Item(3, 9).Value = 0
$ws.Cells.Item(3, 10).Value = 0
$ws.Cells.Item(3, 11).Value = 0
$ws.Cells.Item(3, 12).Value = 0

$ws.Cells.Item(4, 2).Value = 193
$ws.Cells.Item(4, 3).Value = 0
$ws.Cells.Item(4, 4).Value = 41
$ws.Cells.Item(4, 5).Value = 41
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 0
$ws.Cells.Item(4, 8).Value = 0
$ws.Cells.Item(4, 9).Value = 5366
$ws.Cells.Item(4, 10).Value = 0
$ws.Cells.Item(4, 11).Value = 0
$ws.Cells.Item(4, 12).Value = 0

$ws.Cells.Item(5, 2).Value = 1561
$ws.Cells.Item(5, 3).Value = 0
$ws.Cells.Item(5, 4).Value = 59
$ws.Cells.Item(5, 5).Value = 62
$ws.Cells.Item(5, 6).Value = 3
$ws.Cells.Item(5, 7).Value = 0
$ws.Cells.Item(5, 8).Value = 1364
$ws.Cells.Item(5, 9).Value = 28525
$ws.Cells.Item(5, 10).Value = 2
$ws.Cells.Item(5, 11).Value = 0
$ws.Cells.Item(5, 12).Value = 0

$ws.Cells.Item(6, 2).Value = 1912
$ws.Cells.Item(6, 3).Value = 0
$ws.Cells.Item(6, 4).Value = 161
$ws.Cells.Item(6, 5).Value = 172
$ws.Cells.Item(6, 6).Value = 9
$ws.Cells.Item(6, 7).Value = 2
$ws.Cells.Item(6, 8).Value = 2362
$ws.Cells.Item(6, 9).Value = 55549
$ws.Cells.Item(6, 10).Value = 3
$ws.Cells.Item(6, 11).Value = 0
$ws.Cells.Item(6, 12).Value = 0

$ws.Cells.Item(7, 2).Value = 142
$ws.Cells.Item(7, 3).Value = 0
$ws.Cells.Item(7, 4).Value = 17
$ws.Cells.Item(7, 5).Value = 19
$ws.Cells.Item(7, 6).Value = 2
$ws.Cells.Item(7, 7).Value = 0
$ws.Cells.Item(7, 8).Value = 1053
$ws.Cells.Item(7, 9).Value = 7059
$ws.Cells.Item(7, 10).Value = 0
$ws.Cells.Item(7, 11).Value = 0
$ws.Cells.Item(7, 12).Value = 0

$ws.Cells.Item(8, 2).Value = 812
$ws.Cells.Item(8, 3).Value = 1
$ws.Cells.Item(8, 4).Value = 390
$ws.Cells.Item(8, 5).Value = 404
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 2
$ws.Cells.Item(8, 8).Value = 39
$ws.Cells.Item(8, 9).Value = 17307
$ws.Cells.Item(8, 10).Value = 13
$ws.Cells.Item(8, 11).Value = 10
$ws.Cells.Item(8, 12).Value = 386

$ws.Cells.Item(9, 2).Value = 270
$ws.Cells.Item(9, 3).Value = 0
$ws.Cells.Item(9, 4).Value = 10
$ws.Cells.Item(9, 5).Value = 10
$ws.Cells.Item(9, 6).Value = 0
$ws.Cells.Item(9, 7).Value = 0
$ws.Cells.Item(9, 8).Value = 0
$ws.Cells.Item(9, 9).Value = 8000
$ws.Cells.Item(9, 10).Value = 0
$ws.Cells.Item(9, 11).Value = 0
$ws.Cells.Item(9, 12).Value = 0

$ws.Cells.Item(10, 2).Value = 265
$ws.Cells.Item(10, 3).Value = 14
$ws.Cells.Item(10, 4).Value = 277
$ws.Cells.Item(10, 5).Value = 333
$ws.Cells.Item(10, 6).Value = 34
$ws.Cells.Item(10, 7).Value = 8
$ws.Cells.Item(10, 8).Value = 5147
$ws.Cells.Item(10, 9).Value = 10353
$ws.Cells.Item(10, 10).Value = 6
$ws.Cells.Item(10, 11).Value = 0
$ws.Cells.Item(10, 12).Value = 0

$ws.Cells.Item(11, 2).Value = 1057
$ws.Cells.Item(11, 3).Value = 0
$ws.Cells.Item(11, 4).Value = 100
$ws.Cells.Item(11, 5).Value = 101
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 0
$ws.Cells.Item(11, 8).Value = 127
$ws.Cells.Item(11, 9).Value = 17878
$ws.Cells.Item(11, 10).Value = 0
$ws.Cells.Item(11, 11).Value = 0
$ws.Cells.Item(11, 12).Value = 0

$ws.Cells.Item(12, 2).Value = 1407
$ws.Cells.Item(12, 3).Value = 0
$ws.Cells.Item(12, 4).Value = 89
$ws.Cells.Item(12, 5).Value = 113
$ws.Cells.Item(12, 6).Value = 21
$ws.Cells.Item(12, 7).Value = 3
$ws.Cells.Item(12, 8).Value = 2059
$ws.Cells.Item(12, 9).Value = 4336
$ws.Cells.Item(12, 10).Value = 0
$ws.Cells.Item(12, 11).Value = 0
$ws.Cells.Item(12, 12).Value = 0

$ws.Cells.Item(13, 2).Value = 876
$ws.Cells.Item(13, 3).Value = 0
$ws.Cells.Item(13, 4).Value = 185
$ws.Cells.Item(13, 5).Value = 185
$ws.Cells.Item(13, 6).Value = 0
$ws.Cells.Item(13, 7).Value = 0
$ws.Cells.Item(13, 8).Value = 0
$ws.Cells.Item(13, 9).Value = 16330
$ws.Cells.Item(13, 10).Value = 7
$ws.Cells.Item(13, 11).Value = 0
$ws.Cells.Item(13, 12).Value = 0

$ws.Cells.Item(14, 2).Value = 272
$ws.Cells.Item(14, 3).Value = 0
$ws.Cells.Item(14, 4).Value = 25
$ws.Cells.Item(14, 5).Value = 25
$ws.Cells.Item(14, 6).Value = 0
$ws.Cells.Item(14, 7).Value = 0
$ws.Cells.Item(14, 8).Value = 0
$ws.Cells.Item(14, 9).Value = 8000
$ws.Cells.Item(14, 10).Value = 1
$ws.Cells.Item(14, 11).Value = 0
$ws.Cells.Item(14, 12).Value = 0

$ws.Cells.Item(15, 2).Value = 507
$ws.Cells.Item(15, 3).Value = 0
$ws.Cells.Item(15, 4).Value = 38
$ws.Cells.Item(15, 5).Value = 38
$ws.Cells.Item(15, 6).Value = 0
$ws.Cells.Item(15, 7).Value = 0
$ws.Cells.Item(15, 8).Value = 0
$ws.Cells.Item(15, 9).Value = 20238
$ws.Cells.Item(15, 10).Value = 2
$ws.Cells.Item(15, 11).Value = 0
$ws.Cells.Item(15, 12).Value = 0

$ws.Cells.Item(16, 2).Value = 196
$ws.Cells.Item(16, 3).Value = 0
$ws.Cells.Item(16, 4).Value = 1
$ws.Cells.Item(16, 5).Value = 1
$ws.Cells.Item(16, 6).Value = 0
$ws.Cells.Item(16, 7).Value = 0
$ws.Cells.Item(16, 8).Value = 0
$ws.Cells.Item(16, 9).Value = 10000
$ws.Cells.Item(16, 10).Value = 0
$ws.Cells.Item(16, 11).Value = 0
$ws.Cells.Item(16, 12).Value = 0

$ws.Cells.Item(17, 2).Value = 1827
$ws.Cells.Item(17, 3).Value = 0
$ws.Cells.Item(17, 4).Value = 25
$ws.Cells.Item(17, 5).Value = 25
$ws.Cells.Item(17, 6).Value = 0
$ws.Cells.Item(17, 7).Value = 0
$ws.Cells.Item(17, 8).Value = 0
$ws.Cells.Item(17, 9).Value = 18036
$ws.Cells.Item(17, 10).Value = 1
$ws.Cells.Item(17, 11).Value = 0
$ws.Cells.Item(17, 12).Value = 0

$ws.Cells.Item(18, 2).Value = 3903
$ws.Cells.Item(18, 3).Value = 3
$ws.Cells.Item(18, 4).Value = 75
$ws.Cells.Item(18, 5).Value = 180
$ws.Cells.Item(18, 6).Value = 4
$ws.Cells.Item(18, 7).Value = 0
$ws.Cells.Item(18, 8).Value = 1096
$ws.Cells.Item(18, 9).Value = 40465
$ws.Cells.Item(18, 10).Value = 1
$ws.Cells.Item(18, 11).Value = 95
$ws.Cells.Item(18, 12).Value = 14656

